$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D4"  = -7.826700000000001
    "C10" = -13.72049999999999
    "C12" = -13.0023
    "D12" = -8.749899999999997
    "D17" = -8.298399999999994
    "C18" = -11.57119999999999
    "D26" = -7.397000000000005
    "D27" = -7.990699999999999
    "D28" = -7.981
    "C37" = -12.55070000000001
    "D37" = -7.327300000000003
    "C55" = -13.4038
    "D65" = -8.030200000000001
    "C68" = -11.2037
    "D73" = -8.4336
    "C77" = -12.53350000000001
    "C78" = -12.2861
    "D84" = -7.585299999999997
    "D85" = -8.628699999999995
    "D93" = -6.68189999999999
    "D95" = -7.555000000000003
    "D98" = -7.437500000000007
    "D99" = -8.054300000000005
    "D101" = -7.6178
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
